$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update day 21 (row 22) total_venda value
$ws.Range("B22").Value = 25994.16

# Insert a new row for day 22 of 05/2025, shifting the remaining rows down
$ws.Rows(23).Insert()

$ws.Range("A23").Value = 22
$ws.Range("B23").Value = 27147.29
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 2025
$ws.Range("E23").Value = "05/2025"
